# Commit: "Added why and how to initial concept"
# The commit expands the first bullet of the speaker notes on slide 2
# (the notes placeholder) from the short "Problems:" line into a full
# sentence that explains the why/how behind the chosen problems, while
# leaving every other notes bullet untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Locate the notes "body" placeholder (index 2 -> the bullet list, the
# slide-image placeholder is index 1) on the notes page for slide 2.
$notes = $s.NotesPage.Shapes.Placeholders.Item(2)
$tr = $notes.TextFrame.TextRange

# This COM shim always rewrites the whole text frame when TextRange.Text
# is assigned (there is no way to scope the assignment to a single
# paragraph/run), so rebuild the complete bullet list here, only
# changing the very first line's wording - everything else keeps its
# original text exactly as it was.
$newFirstLine = "Problems: The problems we chose were based on the SDG" + [char]0x2019 + "s 9/11 and our aim was to present our research of the innerworkings of Autonomous systems and Digital twin together with the SDG" + [char]0x2019 + "s"

$lines = @(
    $newFirstLine,
    "Too stuffy",
    "Arrows pointing everywhere",
    "Some places miss arrows",
    "Too long nodes",
    "Digital Twins mentioned twice",
    "No clear starting point"
)

$tr.Text = ($lines -join "`n") + "`n"
